$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header of column AC (29) from "nota_iniciativa" to "s_i"
$ws.Cells.Item(1, 29).Value = "s_i"

# Add new column AD (30) header "c4", copying the header formatting (bold + borders)
# from the adjacent AC1 header cell so it matches the rest of row 1.
$ws.Cells.Item(1, 29).Copy() | Out-Null
$ws.Cells.Item(1, 30).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 30).Value = "c4"

# Determine the last data row (97 in the original sheet)
$lastRow = $ws.UsedRange.Rows.Count

# Fill column AD with 0 for every data row, matching the pattern of column AC
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 0
}

$ws.Range("A1").Select() | Out-Null
